$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "'304.17"
$ws.Cells.Item(2,5).Value = "'-0.68%"

$ws.Cells.Item(3,4).Value = "'35.84"
$ws.Cells.Item(3,5).Value = "'0.11%"

$ws.Cells.Item(4,4).Value = "'5.040"
$ws.Cells.Item(4,5).Value = "'-0.55%"

$ws.Cells.Item(5,4).Value = "'0.08018"
$ws.Cells.Item(5,5).Value = "'-0.78%"

$ws.Cells.Item(6,4).Value = "'1.864"
$ws.Cells.Item(6,5).Value = "'-3.73%"

$ws.Cells.Item(7,2).Value = "GateToken"
$ws.Cells.Item(7,3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(7,4).Value = "'4.123"
$ws.Cells.Item(7,5).Value = "'-1.04%"

$ws.Cells.Item(8,2).Value = "KuCoinToken"
$ws.Cells.Item(8,3).Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Cells.Item(8,4).Value = "'7.770"
$ws.Cells.Item(8,5).Value = "'-0.67%"

$ws.Cells.Item(9,2).Value = "MXToken"
$ws.Cells.Item(9,3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(9,4).Value = "'0.9262"
$ws.Cells.Item(9,5).Value = "'-1.61%"

$ws.Cells.Item(10,2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(10,3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(10,4).Value = "'0.1290"
$ws.Cells.Item(10,5).Value = "'-6.10%"

$ws.Cells.Item(11,2).Value = "WazirX"
$ws.Cells.Item(11,3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(11,4).Value = "'0.1898"
$ws.Cells.Item(11,5).Value = "'-0.04%"

$ws.Cells.Item(12,2).Value = "MandalaExchangeToken"
$ws.Cells.Item(12,3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(12,4).Value = "'0.09043"
$ws.Cells.Item(12,5).Value = "'-1.55%"

$ws.Cells.Item(13,2).Value = "BitrueCoin"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(13,4).Value = "'0.03441"
$ws.Cells.Item(13,5).Value = "'-2.09%"

$ws.Cells.Item(14,2).Value = "BitMartToken"
$ws.Cells.Item(14,3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(14,4).Value = "'0.09887"
$ws.Cells.Item(14,5).Value = "'-0.07%"

$ws.Cells.Item(15,2).Value = "BitForexToken"
$ws.Cells.Item(15,3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(15,4).Value = "'0.001401"
$ws.Cells.Item(15,5).Value = "'-3.03%"

$ws.Cells.Item(16,2).Value = "TigerCash"
$ws.Cells.Item(16,3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(16,4).Value = "'0.006292"
$ws.Cells.Item(16,5).Value = "'-5.68%"

$ws.Cells.Item(17,2).Value = "LEO"
$ws.Cells.Item(17,3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(17,4).Value = "'3.838"
$ws.Cells.Item(17,5).Value = "'5.82%"

$ws.Cells.Item(18,4).Value = "'3.375"
$ws.Cells.Item(18,5).Value = "'12.04%"

$ws.Cells.Item(19,4).Value = "'0.3407"
$ws.Cells.Item(19,5).Value = "'-0.49%"

$ws.Cells.Item(20,4).Value = "'0.1334"
$ws.Cells.Item(20,5).Value = "'-0.81%"

$ws.Cells.Item(21,4).Value = "'4.832"
$ws.Cells.Item(21,5).Value = "'-7.90%"

$ws.Cells.Item(22,5).Value = "'-4.94%"

$ws.Cells.Item(23,4).Value = "'0.04366"
$ws.Cells.Item(23,5).Value = "'-0.98%"

$ws.Cells.Item(24,4).Value = "'0.001230"
$ws.Cells.Item(24,5).Value = "'-0.57%"

$ws.Cells.Item(25,4).Value = "'0.004847"
$ws.Cells.Item(25,5).Value = "'1.97%"

$ws.Cells.Item(27,4).Value = "'0.0001301"
$ws.Cells.Item(27,5).Value = "'-0.38%"

$ws.Cells.Item(28,5).Value = "'42.06%"

$ws.Cells.Item(39,4).Value = "'0.01972"
$ws.Cells.Item(39,5).Value = "'-2.65%"

$ws.Cells.Item(40,4).Value = "'0.05147"
$ws.Cells.Item(40,5).Value = "'-0.02%"

$ws.Cells.Item(41,4).Value = "'0.007525"
$ws.Cells.Item(41,5).Value = "'-1.64%"

$ws.Cells.Item(42,5).Value = "'-9.04%"

$ws.Cells.Item(43,4).Value = "'0.1359"
$ws.Cells.Item(43,5).Value = "'-1.10%"

$ws.Cells.Item(44,4).Value = "'0.002112"
$ws.Cells.Item(44,5).Value = "'0.10%"

$ws.Cells.Item(45,4).Value = "'0.009868"
$ws.Cells.Item(45,5).Value = "'-12.72%"

$ws.Cells.Item(46,4).Value = "'0.00006200"
$ws.Cells.Item(46,5).Value = "'-2.26%"

$ws.Cells.Item(47,5).Value = "'-0.56%"

$ws.Cells.Item(48,4).Value = "'64.85"
$ws.Cells.Item(48,5).Value = "'-0.17%"

$ws.Cells.Item(49,4).Value = "'0.001249"
$ws.Cells.Item(49,5).Value = "'-22.00%"

$ws.Cells.Item(50,4).Value = "'0.00002098"
$ws.Cells.Item(50,5).Value = "'-0.56%"

$ws.Cells.Item(51,4).Value = "'0.0001998"
$ws.Cells.Item(51,5).Value = "'-0.56%"
